$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" everywhere it appears ---
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "Handed back: in sync with en-US"

# --- zh-cn: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$zh.Range("I2").Value = "bc7ad01a-d8c4-4593-bace-17fb2811f112.md"
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/203d4f62eaec768c675b42d1c701148cc6893d7a/e2e/bc7ad01a-d8c4-4593-bace-17fb2811f112.md", $null, $null, "bc7ad01a-d8c4-4593-bace-17fb2811f112.md")
$zh.Range("I2").Style = "HyperLink"
$zh.Range("J2").Value = "bc7ad01a-d8c4-4593-bace-17fb2811f112.531a9394bd1e7a4793c0429ba8ba9aa0cc169170.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-13 19:14:40"

# --- de-de: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$de.Range("I2").Value = "bc7ad01a-d8c4-4593-bace-17fb2811f112.md"
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/203d4f62eaec768c675b42d1c701148cc6893d7a/e2e/bc7ad01a-d8c4-4593-bace-17fb2811f112.md", $null, $null, "bc7ad01a-d8c4-4593-bace-17fb2811f112.md")
$de.Range("I2").Style = "HyperLink"
$de.Range("J2").Value = "bc7ad01a-d8c4-4593-bace-17fb2811f112.531a9394bd1e7a4793c0429ba8ba9aa0cc169170.de-de.xlf"
$de.Range("K2").Value = "2016-08-13 19:14:50"

# --- Column width adjustments ---
$overview.Columns.Item(5).ColumnWidth = 29.166666666666664
$overview.Columns.Item(6).ColumnWidth = 29.166666666666664

$zh.Columns.Item(3).ColumnWidth = 29.166666666666664
$zh.Columns.Item(9).ColumnWidth = 39.16666666666667
$zh.Columns.Item(10).ColumnWidth = 39.16666666666667

$de.Columns.Item(3).ColumnWidth = 29.166666666666664
$de.Columns.Item(9).ColumnWidth = 39.16666666666667
$de.Columns.Item(10).ColumnWidth = 39.16666666666667
